$d = $word.ActiveDocument

# Update the date heading (unique text, safe to use Find/Replace)
$d.Content.Find.Execute("2025-05-20 Tuesday", $true, $false, $false, $false, $false, $true, 1, $false, "2025-05-21 Wednesday", 2) | Out-Null

# Update the division-fact table cells by position (some old values repeat,
# so a global Find/Replace would be ambiguous -- address cells directly instead)
$t = $d.Tables.Item(1)

$t.Cell(1, 1).Range.Text = "70÷7=10, 0"
$t.Cell(1, 2).Range.Text = "49÷8=6, 1"
$t.Cell(1, 3).Range.Text = "56÷4=14, 0"
$t.Cell(1, 4).Range.Text = "31÷3=10, 1"
$t.Cell(1, 5).Range.Text = "93÷3=31, 0"
$t.Cell(5, 1).Range.Text = "76÷6=12, 4"
$t.Cell(5, 2).Range.Text = "95÷7=13, 4"
$t.Cell(5, 3).Range.Text = "51÷9=5, 6"
$t.Cell(5, 4).Range.Text = "67÷3=22, 1"
$t.Cell(5, 5).Range.Text = "58÷4=14, 2"
$t.Cell(9, 1).Range.Text = "81÷7=11, 4"
$t.Cell(9, 2).Range.Text = "80÷6=13, 2"
$t.Cell(9, 3).Range.Text = "63÷7=9, 0"
$t.Cell(9, 4).Range.Text = "17÷5=3, 2"
$t.Cell(9, 5).Range.Text = "81÷6=13, 3"
$t.Cell(13, 1).Range.Text = "77÷8=9, 5"
$t.Cell(13, 2).Range.Text = "94÷9=10, 4"
$t.Cell(13, 3).Range.Text = "30÷9=3, 3"
$t.Cell(13, 4).Range.Text = "64÷6=10, 4"
$t.Cell(13, 5).Range.Text = "39÷5=7, 4"
$t.Cell(17, 1).Range.Text = "46÷8=5, 6"
$t.Cell(17, 2).Range.Text = "86÷8=10, 6"
$t.Cell(17, 3).Range.Text = "35÷4=8, 3"
$t.Cell(17, 4).Range.Text = "95÷8=11, 7"
$t.Cell(17, 5).Range.Text = "29÷4=7, 1"
